# Swap the "Integral" theme (currently ppt/theme/theme1.xml, used by the
# slide master / all slides) for the "Office Theme" palette that is already
# embedded in this deck as ppt/theme/theme2.xml (used by the notes master).
#
# The master's ThemeColorScheme exposes the 12 colour slots of the theme's
# <a:clrScheme> in document order:
#   1 dk1  2 lt1  3 dk2  4 lt2  5 accent1 6 accent2
#   7 accent3 8 accent4 9 accent5 10 accent6 11 hlink 12 folHlink
# Writing to ThemeColorScheme.Item(i).RGB rewrites those values directly in
# the theme part backing the slide master (theme1.xml).

function ConvertTo-PptRgb($hex) {
    # PowerPoint's RGB() long is 0x00BBGGRR (i.e. reversed byte order from
    # the usual "RRGGBB" hex notation used in OOXML srgbClr/@val).
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme colour scheme (the palette already present as theme2.xml).
$officeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000",  # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $themeColors.Item($i).RGB = ConvertTo-PptRgb $officeColors[$i - 1]
}
